$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.636.61'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '3.599.35'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '609.61'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.73'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '4.212.40'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '29.85'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').Value = '3.604.63'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').Value = '66.717.18'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.52'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.13'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '427.84'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '3.746.85'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000124'
$ws.Range('E26').Value = '  +5.49%  '
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.39'
$ws.Range('E28').Value = '  +3.11%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.48'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').Value = '3.597.53'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.46'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.87'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '177.71'
$ws.Range('E39').Value = '  +3.96%  '
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.25'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.899'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').Value = '  +9.54%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '25.07'
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '24.24'
$ws.Range('E48').Value = '  +3.96%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').Value = '2.433.17'
$ws.Range('E51').Value = '  +5.54%  '
